$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear T5 and U5 contents
$ws.Range("T5:U5").ClearContents()

# Update the selected cell to I4
$ws.Range("I4").Select()
